# Change com of atoms to mid dist atoms
# Rewrite the residue-combination table (rows 2-21) with the updated
# "mid dist atoms" values, then remove the two now-obsolete trailing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("98, 130, 748, 1073",        2, "130, 130",     "5178, 5582"),
    @("98, 455, 780, 1105",        1, "780",           "5552"),
    @("130, 780, 780",             2, "130, 780",      "3171, 3666"),
    @("98, 130, 748, 780",         1, "130",            "3631"),
    @("130, 455, 748, 1073",       1, "130",            "5131"),
    @("98, 130, 455, 1073",        2, "130, 130",      "4416, 6489"),
    @("423, 1073, 1105, SF, SF",   1, "1105",           "4994"),
    @("130, 455, 748, 780",        1, "130",            "5269"),
    @("455, 748, 1073, 1105",      1, "1105",           "5399"),
    @("423, 748, 780, 1073",       1, "780",            "5677"),
    @("423, 748, 1073, 1105",      2, "1105, 1105",    "5331, 5433"),
    @("423, 455, 748, 1073",       1, "455",            "5886"),
    @("98, 130, 748, 780, 1073",   1, "130",            "6017"),
    @("130, 423, 748, 1073, SF",   1, "130",            "6202"),
    @("98, 130, 455",              1, "130",            "6427"),
    @("130, 455, 780, 1073",       1, "130",            "6562"),
    @("130, 780, 780, 1073",       1, "780",            "6359"),
    @("98, 130, 130, 455, 780",    1, "130",            "6727"),
    @("98, 130, 423, 1073",        1, "130",            "6670"),
    @("98, 98, 455, 455",          1, "455",            "6748")
)

# Columns C and D hold text (e.g. "780", "5178, 5582") even when the text
# looks numeric, so force a text format before assigning the values - this
# keeps Excel from silently re-interpreting them as numbers.
$ws.Range("C2:D21").NumberFormat = "@"

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $row = $row + 1
}

# The new table only spans down to row 21; remove the two trailing rows
# that belonged to the old (longer) table.
$ws.Rows.Item(23).Delete()
$ws.Rows.Item(22).Delete()
